$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4744.7
$ws.Range("I19").Value = 6713.857
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 6713.857
$ws.Range("L19").Value = 150
$ws.Range("M19").Value = -6538.857
$ws.Range("N19").Value = -500
$ws.Range("H100").Value = 2788.7222
$ws.Range("I100").Value = 2046.1538
$ws.Range("K100").Value = 2046.1538
$ws.Range("M100").Value = -1505.1538
$ws.Range("H112").Value = 1786.55
$ws.Range("J112").Value = 1874.2354
$ws.Range("L112").Value = 5622.706200000001
$ws.Range("N112").Value = -7838.706200000001
$ws.Range("H127").Value = 848.3333
$ws.Range("J127").Value = 555
$ws.Range("L127").Value = 1665
$ws.Range("N127").Value = -11585
$ws.Range("H129").Value = 1885.3684
$ws.Range("I129").Value = 713.9
$ws.Range("K129").Value = 2141.7
$ws.Range("M129").Value = 2858.3
$ws.Range("H132").Value = 2422.5
$ws.Range("I132").Value = 1220.4651
$ws.Range("K132").Value = 3661.3953
$ws.Range("M132").Value = -1131.3953

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 14999.75
$ws.Range("J44").Value = 14999.75
$ws.Range("L44").Value = 14999.75
$ws.Range("N44").Value = -15975.75
$ws.Range("H45").Value = 10047.125
$ws.Range("I45").Value = 15716.125
$ws.Range("J45").Value = 4378.125
$ws.Range("K45").Value = 15716.125
$ws.Range("L45").Value = 4378.125
$ws.Range("M45").Value = -15339.125
$ws.Range("N45").Value = -5132.125
$ws.Range("H102").Value = 2677.1924
$ws.Range("I102").Value = 2737.7917
$ws.Range("J102").Value = 1950
$ws.Range("K102").Value = 2737.7917
$ws.Range("L102").Value = 1950
$ws.Range("M102").Value = -1115.7917
$ws.Range("N102").Value = -5194
$ws.Range("H110").Value = 1355.4166
$ws.Range("I110").Value = 1355.4166
$ws.Range("K110").Value = 1355.4166
$ws.Range("M110").Value = 689.5834

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4333.3335
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("H107").Value = 102789.1
$ws.Range("I107").Value = 144413.14
$ws.Range("J107").Value = 5666.3335
$ws.Range("K107").Value = 144413.14
$ws.Range("L107").Value = 5666.3335
$ws.Range("M107").Value = -142493.14
$ws.Range("N107").Value = -9506.333500000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14861
$ws.Range("I31").Value = 4126.5386
$ws.Range("J31").Value = 22613.666
$ws.Range("K31").Value = 4126.5386
$ws.Range("L31").Value = 22613.666
$ws.Range("M31").Value = -3831.5386
$ws.Range("N31").Value = -23203.666
$ws.Range("H34").Value = 14861
$ws.Range("I34").Value = 4126.5386
$ws.Range("J34").Value = 22613.666
$ws.Range("K34").Value = 4126.5386
$ws.Range("L34").Value = 22613.666
$ws.Range("M34").Value = -3924.5386
$ws.Range("N34").Value = -23017.666
$ws.Range("H122").Value = 70150.8
$ws.Range("I122").Value = 91620
$ws.Range("K122").Value = 274860
$ws.Range("M122").Value = -272410
$ws.Range("H132").Value = 2166.5334
$ws.Range("I132").Value = 2166.5334
$ws.Range("K132").Value = 6499.600199999999
$ws.Range("M132").Value = -3969.600199999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2715
$ws.Range("J39").Value = 3947.5
$ws.Range("L39").Value = 11842.5
$ws.Range("N39").Value = -12430.5
$ws.Range("H50").Value = 1500307
$ws.Range("I50").Value = 258
$ws.Range("K50").Value = 774
$ws.Range("M50").Value = -293
$ws.Range("H53").Value = 1500307
$ws.Range("I53").Value = 258
$ws.Range("K53").Value = 774
$ws.Range("M53").Value = -293
$ws.Range("H106").Value = 5682.7617
$ws.Range("J106").Value = 5766.95
$ws.Range("L106").Value = 17300.85
$ws.Range("N106").Value = -19192.85
$ws.Range("H131").Value = 1876.4
$ws.Range("I131").Value = 2055
$ws.Range("J131").Value = 1757.3334
$ws.Range("K131").Value = 6165
$ws.Range("L131").Value = 5272.0002
$ws.Range("M131").Value = -1125
$ws.Range("N131").Value = -15352.0002
$ws.Range("H132").Value = 1872.0869
$ws.Range("I132").Value = 996.7778
$ws.Range("K132").Value = 8971.0002
$ws.Range("M132").Value = -6441.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6845.4546
$ws.Range("I70").Value = 6705.5557
$ws.Range("J70").Value = 7475
$ws.Range("K70").Value = 6705.5557
$ws.Range("L70").Value = 7475
$ws.Range("M70").Value = -6435.5557
$ws.Range("N70").Value = -8015
$ws.Range("H73").Value = 6845.4546
$ws.Range("I73").Value = 6705.5557
$ws.Range("J73").Value = 7475
$ws.Range("K73").Value = 6705.5557
$ws.Range("L73").Value = 7475
$ws.Range("M73").Value = -5769.5557
$ws.Range("N73").Value = -9347
$ws.Range("H80").Value = 3836.625
$ws.Range("I80").Value = 2699.75
$ws.Range("J80").Value = 4973.5
$ws.Range("K80").Value = 2699.75
$ws.Range("L80").Value = 4973.5
$ws.Range("M80").Value = -1701.75
$ws.Range("N80").Value = -6969.5
$ws.Range("H83").Value = 3836.625
$ws.Range("I83").Value = 2699.75
$ws.Range("J83").Value = 4973.5
$ws.Range("K83").Value = 13498.75
$ws.Range("L83").Value = 24867.5
$ws.Range("M83").Value = -8506.75
$ws.Range("N83").Value = -34851.5
$ws.Range("H132").Value = 3363.7058
$ws.Range("I132").Value = 3805.3076
$ws.Range("J132").Value = 1928.5
$ws.Range("K132").Value = 11415.9228
$ws.Range("L132").Value = 5785.5
$ws.Range("M132").Value = -8885.9228
$ws.Range("N132").Value = -10845.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 101770
$ws.Range("I16").Value = 1528.2858
$ws.Range("K16").Value = 1528.2858
$ws.Range("M16").Value = -1358.2858
$ws.Range("H46").Value = 25079.3
$ws.Range("I46").Value = 64727.43
$ws.Range("J46").Value = 3730.3076
$ws.Range("K46").Value = 64727.43
$ws.Range("L46").Value = 3730.3076
$ws.Range("M46").Value = -64539.43
$ws.Range("N46").Value = -4106.3076
$ws.Range("H55").Value = 463.91666
$ws.Range("I55").Value = 561.5
$ws.Range("J55").Value = 366.33334
$ws.Range("K55").Value = 561.5
$ws.Range("L55").Value = 366.33334
$ws.Range("M55").Value = -388.5
$ws.Range("N55").Value = -712.33334
$ws.Range("H68").Value = 3766.4
$ws.Range("I68").Value = 3742.3333
$ws.Range("J68").Value = 3802.5
$ws.Range("K68").Value = 3742.3333
$ws.Range("L68").Value = 3802.5
$ws.Range("M68").Value = -2993.3333
$ws.Range("N68").Value = -5300.5
$ws.Range("H71").Value = 3766.4
$ws.Range("I71").Value = 3742.3333
$ws.Range("J71").Value = 3802.5
$ws.Range("K71").Value = 18711.6665
$ws.Range("L71").Value = 19012.5
$ws.Range("M71").Value = -14967.6665
$ws.Range("N71").Value = -26500.5
$ws.Range("H82").Value = 1281.92
$ws.Range("I82").Value = 1153.2354
$ws.Range("K82").Value = 1153.2354
$ws.Range("M82").Value = -792.2354
$ws.Range("H85").Value = 1281.92
$ws.Range("I85").Value = 1153.2354
$ws.Range("K85").Value = 1153.2354
$ws.Range("M85").Value = 94.76459999999997

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 3375.8235
$ws.Range("I100").Value = 3256.5
$ws.Range("K100").Value = 6513
$ws.Range("M100").Value = -5972
$ws.Range("H132").Value = 1943.2
$ws.Range("I132").Value = 1430.3158
$ws.Range("J132").Value = 3567.3333
$ws.Range("K132").Value = 4290.9474
$ws.Range("L132").Value = 10701.9999
$ws.Range("M132").Value = -1760.9474
$ws.Range("N132").Value = -15761.9999
$ws.Range("H136").Value = 3069.4727
$ws.Range("I136").Value = 2696.1333
$ws.Range("J136").Value = 4749.5
$ws.Range("K136").Value = 8088.3999
$ws.Range("L136").Value = 14248.5
$ws.Range("M136").Value = -5538.3999
$ws.Range("N136").Value = -19348.5

Write-Output "Applied $([int]202) cell updates across 8 sheets"
